# Adds the latest Bundesliga matchday (Matchday 6, 2025-10) results to the
# end of the match table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newMatches = @(
    @(540458, "2025-10-03T18:30:00Z", 6,  2, "TSG 1899 Hoffenheim",      1, "1. FC Köln",             0, 1, "AwayWin"),
    @(540451, "2025-10-04T13:30:00Z", 6,  3, "Bayer 04 Leverkusen",     28, "1. FC Union Berlin",     2, 0, "HomeWin"),
    @(540453, "2025-10-04T13:30:00Z", 6,  4, "Borussia Dortmund",      721, "RB Leipzig",             1, 1, "Draw"),
    @(540454, "2025-10-04T13:30:00Z", 6, 12, "SV Werder Bremen",        20, "FC St. Pauli 1910",      1, 0, "HomeWin"),
    @(540457, "2025-10-04T13:30:00Z", 6, 16, "FC Augsburg",             11, "VfL Wolfsburg",          3, 1, "HomeWin"),
    @(540452, "2025-10-04T16:30:00Z", 6, 19, "Eintracht Frankfurt",      5, "FC Bayern München",      0, 3, "AwayWin"),
    @(540455, "2025-10-05T13:30:00Z", 6, 10, "VfB Stuttgart",           44, "1. FC Heidenheim 1846",  1, 0, "HomeWin"),
    @(540459, "2025-10-05T15:30:00Z", 6,  7, "Hamburger SV",            15, "1. FSV Mainz 05",        4, 0, "HomeWin"),
    @(540456, "2025-10-05T17:30:00Z", 6, 18, "Borussia Mönchengladbach", 17, "SC Freiburg",           0, 0, "Draw")
)

$startRow = $ws.UsedRange.Rows.Count + 1

for ($i = 0; $i -lt $newMatches.Count; $i++) {
    $row = $startRow + $i
    $data = $newMatches[$i]
    for ($col = 1; $col -le $data.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$col - 1]
    }
}
